$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "7.08") need to be forced
# to stay text: flip to a text format before the write, then restore the default
# "Normal" style afterwards so no style index is left behind on the cell.

$ws.Range("D2").Value = "67.287.13"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.481.70"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "3.483.76"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +4.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "4.082.27"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.77%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "67.323.76"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "3.481.28"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.871"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.89%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "2.846.81"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "336.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -2.29%  "
